$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 2004年-2009年 rows (old rows 2-7); the remaining data (2010年-2020年)
# shifts up to rows 2-12, matching the target layout (dimension becomes A1:H14
# once the two new rows below are appended).
$ws.Rows("2:7").Delete()

# Append the new 2021年 row (row 13), copying row 12's formatting for column A
# (bold/bordered/centered header-style cell used throughout column A).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = 19
$ws.Range("G13").Value = 9
$ws.Range("H13").Value = 1065200
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Append the new 2022年 row (row 14) -- only the "地震灾害次数" count is known.
$ws.Range("A14").Value = "2022年"
$ws.Range("F14").Value = 27
$ws.Range("A12").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$excel.CutCopyMode = $false
